# Search user using username
# Adds three new rows of locator data to the "Web" worksheet:
#   UM_searchUsername, UM_searchButton, UM_searchResultUsername
# They fill previously-blank rows 16-18 (A:C) and a brand-new row 21 is
# appended as a blank row at the bottom (matching the source diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Web")

# Row 16: search username field
$ws.Range("A16").Value = "UM_searchUsername"
$ws.Range("B16").Value = "//*[@id=`"app`"]/div[1]/div[2]/div[2]/div/div[1]/div[2]/form/div[1]/div/div[1]/div/div[2]/input"
$ws.Range("C16").Value = "By.xpath"

# Row 17: search button
$ws.Range("A17").Value = "UM_searchButton"
$ws.Range("B17").Value = "(//button[normalize-space()='Search'])"
$ws.Range("C17").Value = "By.xpath"

# Row 18: search result username
$ws.Range("A18").Value = "UM_searchResultUsername"
$ws.Range("B18").Value = "(//div[contains(text(),'Admin')])"
$ws.Range("C18").Value = "By.xpath"

# Row 21: new blank row appended at the bottom (empty strings).
# A bare "'" collapses to an empty string (same as the existing blank
# rows 19-20) while still materializing the cell/row and extending the
# sheet's used range/dimension to A1:C21, which a plain "" assignment
# does not do (the engine treats "" as "clear cell", removing it again).
$ws.Range("A21").Value = "'"
$ws.Range("B21").Value = "'"
$ws.Range("C21").Value = "'"
